$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily sales record (day 11, June/2025) was added to the data,
# inserted right before the existing "May 2025" block (which previously
# started at row 12). Inserting a row at 12 shifts all the following
# rows (old rows 12-103, the May/April/March data) down by one, to rows
# 13-104, matching the diff.
$ws.Rows("12:12").Insert()

$ws.Cells.Item(12, 1).Value2 = 11
$ws.Cells.Item(12, 2).Value2 = 33119.36
$ws.Cells.Item(12, 3).Value2 = 6
$ws.Cells.Item(12, 4).Value2 = 2025
$ws.Cells.Item(12, 5).Value2 = "06/2025"
